# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New computed "K" values (column G) for rows 2-29
$kValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 1
    6  = 2
    7  = 0
    8  = 2
    9  = 2
    10 = 0
    11 = 1
    12 = 0
    13 = 2
    14 = 2
    15 = 1
    16 = 0
    17 = 2
    18 = 2
    19 = 2
    20 = 1
    21 = 2
    22 = 1
    23 = 2
    24 = 1
    25 = 0
    26 = 0
    27 = 1
    28 = 0
    29 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
